# Re-order the codeforiati: group/category columns.
#
# Before:  D = codeforiati:group-code     E = codeforiati:category-name
#          F = codeforiati:category-code  G = codeforiati:group-name
#
# After:   D = codeforiati:category-name  E = codeforiati:group-name
#          F = codeforiati:category-code  G = codeforiati:group-code
#
# i.e. for every row: new D = old E, new E = old G, new G = old D, F stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Header row (row 1)
$ws.Cells.Item(1, 4).Value2 = "codeforiati:category-name"
$ws.Cells.Item(1, 5).Value2 = "codeforiati:group-name"
$ws.Cells.Item(1, 6).Value2 = "codeforiati:category-code"
$ws.Cells.Item(1, 7).Value2 = "codeforiati:group-code"

for ($r = 2; $r -le $lastRow; $r++) {
    $oldD = $ws.Cells.Item($r, 4).Value2
    $oldE = $ws.Cells.Item($r, 5).Value2
    $oldG = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 4).Value2 = $oldE
    $ws.Cells.Item($r, 5).Value2 = $oldG
    $ws.Cells.Item($r, 7).Value2 = $oldD
}
